# Insert 4 new daily price rows above the current row 759 (pushing the
# existing 759-806 block down to 763-810) and populate them with the new
# "Cebolla" (onion) price observations for Vega Modelo de Temuco.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 4 blank rows at 759..762; everything currently at 759-806 shifts
# down to 763-810 (dimension grows from R806 to R810 automatically).
$ws.Rows("759:762").Insert()

# Common columns shared by every row in this block.
$mercadoId = 10
$mercado = "Vega Modelo de Temuco"
$region = "La Araucanía"
$codreg = 9
$categoriaId = 100112004
$categoria = "Cebolla"
$clasificacion = "Hortaliza"

$newRows = @(
    @{ Row = 759; Fecha = 44585; Variedad = "Morada(o)";        Calidad = "1a (guarda)"; Volumen = 155; PMin = 10000; PMax = 10000; PProm = 10000; Unidad = "`$/malla 18 kilos"; Origen = "Región de O'Higgins"; PKg = 556; KgUnidades = 18 },
    @{ Row = 760; Fecha = 44585; Variedad = "Sin especificar";  Calidad = "1a nueva(o)"; Volumen = 450; PMin = 5000;  PMax = 5000;  PProm = 5000;  Unidad = "`$/malla 18 kilos"; Origen = "Región de O'Higgins"; PKg = 278; KgUnidades = 18 },
    @{ Row = 761; Fecha = 44585; Variedad = "Sin especificar";  Calidad = "1a nueva(o)"; Volumen = 450; PMin = 4500;  PMax = 5000;  PProm = 4778;  Unidad = "`$/malla 18 kilos"; Origen = "Región del Maule";      PKg = 265; KgUnidades = 18 },
    @{ Row = 762; Fecha = 44585; Variedad = "Sin especificar";  Calidad = "Primera";     Volumen = 750; PMin = 4500;  PMax = 5000;  PProm = 4767;  Unidad = "`$/malla 18 kilos"; Origen = "Perú";                    PKg = 265; KgUnidades = 18 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $r.Fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $categoriaId
    $ws.Cells.Item($row, 7).Value = $categoria
    $ws.Cells.Item($row, 8).Value = $r.Variedad
    $ws.Cells.Item($row, 9).Value = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.PMin
    $ws.Cells.Item($row, 12).Value = $r.PMax
    $ws.Cells.Item($row, 13).Value = $r.PProm
    $ws.Cells.Item($row, 14).Value = $r.Unidad
    $ws.Cells.Item($row, 15).Value = $r.Origen
    $ws.Cells.Item($row, 16).Value = $r.PKg
    $ws.Cells.Item($row, 17).Value = $r.KgUnidades
    $ws.Cells.Item($row, 18).Value = $clasificacion
}
